# Append two new data rows (73, 74) to Sheet1, mirroring the existing
# row layout (stockname, alert_date, rowNumber, date1, row1, value1,
# date2, row2, value2, buyORsell, slope, intercept, window_size,
# percentage_value, two_line_count).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 73
$ws.Cells.Item(73, 1).Value = "PEL.NS"
$ws.Cells.Item(73, 2).Value = 37803
$ws.Cells.Item(73, 3).Value = 30
$ws.Cells.Item(73, 4).Value = 37073
$ws.Cells.Item(73, 5).Value = 22
$ws.Cells.Item(73, 6).Value = 21.93608474731445
$ws.Cells.Item(73, 7).Value = 37438
$ws.Cells.Item(73, 8).Value = 26
$ws.Cells.Item(73, 9).Value = 21.79438591003418
$ws.Cells.Item(73, 10).Value = "Low"
$ws.Cells.Item(73, 11).Value = -0.03542470932006836
$ws.Cells.Item(73, 12).Value = 22.71542835235596
$ws.Cells.Item(73, 13).Value = 3
$ws.Cells.Item(73, 14).Value = 1
$ws.Cells.Item(73, 15).Value = 2

# Row 74
$ws.Cells.Item(74, 1).Value = "PEL.NS"
$ws.Cells.Item(74, 2).Value = 44835
$ws.Cells.Item(74, 3).Value = 107
$ws.Cells.Item(74, 4).Value = 42826
$ws.Cells.Item(74, 5).Value = 85
$ws.Cells.Item(74, 6).Value = 1741.661987304688
$ws.Cells.Item(74, 7).Value = 44470
$ws.Cells.Item(74, 8).Value = 103
$ws.Cells.Item(74, 9).Value = 1751.69140625
$ws.Cells.Item(74, 10).Value = "High"
$ws.Cells.Item(74, 11).Value = 0.55718994140625
$ws.Cells.Item(74, 12).Value = 1694.300842285156
$ws.Cells.Item(74, 13).Value = 3
$ws.Cells.Item(74, 14).Value = 1
$ws.Cells.Item(74, 15).Value = 2

# Match the date-formatted style used by the other alert_date/date1/date2
# columns (B, D, G) in the existing rows.
$ws.Range("B73").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D73").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G73").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B74").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D74").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G74").NumberFormat = "YYYY-MM-DD HH:MM:SS"
